$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values to reflect repulled data / recalculated mean
$ws.Range("F3").Value = 1
$ws.Range("F6").Value = -3
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = -5
